$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.382.93'
$ws.Range('E2').Value = '  +4.03%  '
$ws.Range('D3').Value = '1.594.42'
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = "'214.86"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.89%  '
$ws.Range('D6').Value = "'0.494"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.92%  '
$ws.Range('E7').Value = '  -0.04%  '
$ws.Range('D8').Value = "'24.12"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +8.54%  '
$ws.Range('E9').Value = '  +0.72%  '
$ws.Range('D10').Value = "'0.0601"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +0.79%  '
$ws.Range('D11').Value = "'0.0888"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.98%  '
$ws.Range('D12').Value = '1.821.26'
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').Value = '1.596.29'
$ws.Range('E13').Value = '  +1.99%  '
$ws.Range('D14').Value = "'0.532"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.38%  '
$ws.Range('E15').Value = '  -0.33%  '
$ws.Range('D16').Value = '28.385.41'
$ws.Range('E16').Value = '  +4.12%  '
$ws.Range('D17').Value = "'63.17"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.86%  '
$ws.Range('D18').Value = "'228.13"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.62%  '
$ws.Range('D19').Value = '0.0₃0711'
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('D20').Value = "'7.50"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.15%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').Value = "'4.10"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.35%  '
$ws.Range('E23').Value = '  -0.59%  '
$ws.Range('D24').Value = "'1.96"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.80%  '
$ws.Range('D25').Value = "'151.69"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.12%  '
$ws.Range('E26').Value = '  +0.93%  '
$ws.Range('E27').Value = '  +0.32%  '
$ws.Range('D28').Value = "'6.59"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.74%  '
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('D31').Value = "'0.0477"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +1.05%  '
$ws.Range('D32').Value = "'3.24"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -0.20%  '
$ws.Range('D33').Value = "'3.14"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.09%  '
$ws.Range('D34').Value = '1.396.78'
$ws.Range('E34').Value = '  -4.34%  '
$ws.Range('D35').Value = "'1.60"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('E36').Value = '  -5.57%  '
$ws.Range('E37').Value = '  +0.44%  '
$ws.Range('E38').Value = '  +0.46%  '
$ws.Range('D39').Value = "'2.51"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +6.95%  '
$ws.Range('D40').Value = "'0.541"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.40%  '
$ws.Range('D41').Value = "'0.815"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').Value = "'5.72"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.80%  '
$ws.Range('E43').Value = '  -0.06%  '
$ws.Range('D44').Value = "'1.88"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +6.92%  '
$ws.Range('D45').Value = "'0.984"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.26%  '
$ws.Range('D46').Value = "'64.47"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.19%  '
$ws.Range('D47').Value = '1.731.59'
$ws.Range('E47').Value = '  +1.77%  '
$ws.Range('D48').Value = "'87.44"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.60%  '
$ws.Range('D49').Value = "'2.13"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.07%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').Value = "'0.0525"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.04%  '
